# Updates from Sadia - 17th May
# Adds 3 new test-case rows (TC_013/014/015 -> Sauce Labs Backpack / Bolt
# T-Shirt / Bike Light), a new "hilo" sort-value entry, and three new
# "checkout information" columns (firstName, lastName, postalCode) filled
# in on the last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells for the checkout-info columns -----------------------
$ws.Range("M1").Value = "firstName"
$ws.Range("N1").Value = "lastName"
$ws.Range("O1").Value = "postalCode"

# --- Row 15: SortValue moved from column L to column K --------------------
$ws.Range("L15").ClearContents()
$ws.Range("K15").Value = "hilo"

# --- Row 16: TC_013 / Sauce Labs Backpack ----------------------------------
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "TC_013"
$ws.Range("C16").Value = "standard_user"
$ws.Range("D16").Value = "secret_sauce"
$ws.Range("L16").Value = "Sauce Labs Backpack"

# --- Row 17: TC_014 / Sauce Labs Bolt T-Shirt ------------------------------
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "TC_014"
$ws.Range("C17").Value = "standard_user"
$ws.Range("D17").Value = "secret_sauce"
$ws.Range("L17").Value = "Sauce Labs Bolt T-Shirt"

# --- Row 18: TC_015 / Sauce Labs Bike Light + checkout info ---------------
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "TC_015"
$ws.Range("C18").Value = "standard_user"
$ws.Range("D18").Value = "secret_sauce"
$ws.Range("L18").Value = "Sauce Labs Bike Light"
$ws.Range("M18").Value = "Sadia"
$ws.Range("N18").Value = "Nasim"
$ws.Range("O18").Value = 711101

# --- Selection + page setup, matching the saved workbook state -----------
$ws.Range("S14").Select()
$ws.PageSetup.Orientation = 1
